$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.984.52"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "1.564.52"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.11"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("E10").Value = "  +2.15%  "

$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("D12").Value = "1.786.65"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").Value = "1.564.21"
$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("E15").Value = "  -0.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").Value = "26.987.24"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "0.0₃0706"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.31%  "

$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("E22").Value = "  +0.88%  "

$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("E28").Value = "  +0.87%  "

$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("E32").Value = "  -0.24%  "

$ws.Range("E33").Value = "  +0.91%  "

$ws.Range("D34").Value = "1.424.45"
$ws.Range("E34").Value = "  -1.47%  "

$ws.Range("E35").Value = "  +2.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.70%  "

$ws.Range("E37").Value = "  +1.91%  "

$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.536"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.00%  "

$ws.Range("E40").Value = "  +0.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "

$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("E43").Value = "  +2.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").Value = "1.701.02"
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0961"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.19%  "
